# Weekly update for "Vega Monumental Concepción - Cebollín": a new
# week's worth of data (two quality lines: Primera / Segunda) is
# prepended above the existing history, so the previously-existing
# rows 14-17 shift down to rows 16-19 unchanged, and the two new rows
# (14-15) carry the newest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 14, pushing the old
# rows 14-17 down to 16-19 (with their values/formatting intact).
$ws.Range("A14:R15").Insert()

# --- Row 14: newest week, quality "Primera" ---
$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 44491
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 100112037
$ws.Cells.Item(14, 7).Value = "Cebollín"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 200
$ws.Cells.Item(14, 11).Value = 600
$ws.Cells.Item(14, 12).Value = 700
$ws.Cells.Item(14, 13).Value = 650
$ws.Cells.Item(14, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(14, 15).Value = "Región Metropolitana"
$ws.Cells.Item(14, 16).Value = 108
$ws.Cells.Item(14, 17).Value = 6
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# --- Row 15: newest week, quality "Segunda" ---
$ws.Cells.Item(15, 1).Value = 11
$ws.Cells.Item(15, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(15, 3).Value = "Bíobío"
$ws.Cells.Item(15, 4).Value = 44491
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 6).Value = 100112037
$ws.Cells.Item(15, 7).Value = "Cebollín"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Segunda"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 500
$ws.Cells.Item(15, 12).Value = 500
$ws.Cells.Item(15, 13).Value = 500
$ws.Cells.Item(15, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 83
$ws.Cells.Item(15, 17).Value = 6
$ws.Cells.Item(15, 18).Value = "Hortaliza"
